$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cell "Other Attributes" spanning D3:I3 -----------------
# Bold the whole header row (B3:I3) first - matches fontId=1 (bold) cellXfs
# entries created for B3/C3 (s=2) and D3:I3 (s=3).
$ws.Range("B3:I3").Font.Bold = $true

$ws.Range("D3").Value = "Other Attributes"
$ws.Range("D3:I3").HorizontalAlignment = -4108   # xlCenter

# Merge the new header cells across D3:I3 (like the existing B-column merges)
$ws.Range("D3:I3").MergeCells = $true

# --- Column D got a touch wider to fit the new content -----------------
$ws.Columns("D").ColumnWidth = 12.6640625

# --- Rename the "list_email_*" SK rows to "list_list_email_*" ----------
# (rows 5/6 under the first user_email_1 GSI block, and rows 8/9 under the
# second user_email_2 GSI block)
$ws.Range("C5").Value = "list_list_email_1"
$ws.Range("C6").Value = "list_list_email_2"
$ws.Range("C8").Value = "list_list_email_1"
$ws.Range("C9").Value = "list_list_email_2"

# --- Restore the active selection to F6 (as captured in the saved file) -
$ws.Range("F6").Select()

Write-Host "done"
